$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The generator (horario.py) was collapsing every "resultado de aprendizaje"
# tied to the same competence/course into a single overwritten row. The fix
# emits one row per learning outcome, leaving the course-code (A) and hour
# count (B) columns blank on the repeated rows (they're only shown once,
# on the row where the block starts).

# --- Block 1: "Negociación, Marketing Digital" (course 2651969 A) -----------
# Row 28 already holds the first outcome. Insert 3 more rows right after it.
$ws.Range("A29:A31").EntireRow.Insert()

$ws.Range("C29").Value = "MULTIMEDIA"
$ws.Range("D29").Value = "Negociación, Marketing Digital"
$ws.Range("E29").Value = "Realizar negociación con los proveedores y clientes, según los objetivos y estrategias establecidas por la organización."
$ws.Range("F29").Value = "Realizar negociación con los proveedores y clientes, según los objetivos y estrategias establecidas por la organización."
$ws.Range("G29").Value = "Crear la propuesta y/o contrapropuesta que establezca detalles de rentabilidad, precios, y compromisos con base en las políticas del cliente (financiación, presupuesto, objetivos e identidad corporativa). (Negociación – Marketing Digital)"

$ws.Range("C30").Value = "MULTIMEDIA"
$ws.Range("D30").Value = "Negociación, Marketing Digital"
$ws.Range("E30").Value = "Realizar negociación con los proveedores y clientes, según los objetivos y estrategias establecidas por la organización."
$ws.Range("F30").Value = "Realizar negociación con los proveedores y clientes, según los objetivos y estrategias establecidas por la organización."
$ws.Range("G30").Value = "Identificar las tendencias del mercado y del diseño en la producción de proyectos multimedia para orientar al cliente. "

$ws.Range("C31").Value = "MULTIMEDIA"
$ws.Range("D31").Value = "Negociación, Marketing Digital"
$ws.Range("E31").Value = "Realizar negociación con los proveedores y clientes, según los objetivos y estrategias establecidas por la organización."
$ws.Range("F31").Value = "Realizar negociación con los proveedores y clientes, según los objetivos y estrategias establecidas por la organización."
$ws.Range("G31").Value = "Definir el costo del proyecto teniendo en cuenta los gastos fijos, variables y tiempo invertido en la realización del mismo. "

# Row 32 is the former row 29 ("Proyecto (pruebas de usuario)" / 2651969 A),
# already pushed down by the insert above - its content stays as-is.

# --- Block 2: "Proyecto (pruebas de usuario)" (course 2651969 A) -----------
# Insert 2 more rows right after row 32.
$ws.Range("A33:A34").EntireRow.Insert()

$ws.Range("C33").Value = "MULTIMEDIA"
$ws.Range("D33").Value = "Proyecto (pruebas de usuario)"
$ws.Range("E33").Value = "Entregar la aplicación multimedia para evaluar la satisfacción del cliente"
$ws.Range("F33").Value = "Entregar la aplicación multimedia para evaluar la satisfacción del cliente"
$ws.Range("G33").Value = "Realizar las modificaciones pertinentes de acuerdo con lo evaluado en las pruebas de accesibilidad, diseño, escalabilidad y usabilidad de la multimedia."

$ws.Range("C34").Value = "MULTIMEDIA"
$ws.Range("D34").Value = "Proyecto (pruebas de usuario)"
$ws.Range("E34").Value = "Entregar la aplicación multimedia para evaluar la satisfacción del cliente"
$ws.Range("F34").Value = "Entregar la aplicación multimedia para evaluar la satisfacción del cliente"
$ws.Range("G34").Value = "Elaborar los manuales y ayudas análogas o digitales necesarias para facilitar la operación del proyecto multimedia."

# Row 35 is the former row 30 (2771153 A / INT CONT DIGITALES), already
# pushed down by the insert above - its content stays as-is.

Write-Output "Inserted 5 rows; sheet now spans A2:G35"
